# Fixing more bugs with select others
# This script edits the "survey" sheet to add a new "inputAttributes.data-type"
# column, add a horizontal select example row, add a content-provider-query
# select_one example row, and updates the "queries" sheet accordingly.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")
$queries = $wb.Worksheets.Item("queries")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# 1. survey sheet: insert a new column D ("inputAttributes.data-type") right
#    before the existing "condition" column. This shifts condition/name/label
#    from D/E/F to E/F/G.
# ---------------------------------------------------------------------------
$survey.Columns.Item(4).Insert()

$survey.Range("D1").Value = "inputAttributes.data-type"

# Restore / set the custom column widths for the affected columns. (E and F
# keep the widths that moved with them automatically; D is brand new and G
# needs a new, wider width since it now holds the "label" column.)
$survey.Columns.Item(4).ColumnWidth = 19.833
$survey.Columns.Item(7).ColumnWidth = 45.166

# ---------------------------------------------------------------------------
# 2. survey sheet: new row 16 - horizontal select appearance example.
# ---------------------------------------------------------------------------
$survey.Range("B16").Value = "select_one yes_no"
$survey.Range("D16").Value = "horizontal"
$survey.Range("F16").Value = "h_select"
$survey.Range("G16").Value = "Horizontal select example."

# ---------------------------------------------------------------------------
# 3. survey sheet: new row 17 - content provider query select_one example.
# ---------------------------------------------------------------------------
$survey.Range("B17").Value = "select_one content_provider_test"
$survey.Range("F17").Value = "cp_test"
$survey.Range("G17").Value = "This demos a content provider query."

# ---------------------------------------------------------------------------
# 4. queries sheet: rename the "odk_values" query to "content_provider_test"
#    and give it a callback function in column C.
# ---------------------------------------------------------------------------
$queries.Range("A5").Value = "content_provider_test"
$queries.Range("C5").Value = '[{ name: "test", label : JSON.stringify(context) }]'
